$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the existing "Date Ordered" row (old row 4),
# pushing everything below down by two rows (formulas auto-adjust).
$ws.Rows("4:5").Insert()

# New row 5 (blank spacer row 4 is left as-is, just bold-styled by the insert):
# "Order Lodged with" / "UC ECE Dept."
$ws.Range("B5").Value = "Order Lodged with"
$ws.Range("B5").Font.Bold = $true
$ws.Range("C5").Value = "UC ECE Dept."

# The "Date Ordered" label (now on row 6) becomes bold too.
$ws.Range("B6").Font.Bold = $true

# Widen column C slightly (stored width becomes 13).
$ws.Columns("C").ColumnWidth = 12.15

# Update the active selection.
$ws.Range("E6").Select() | Out-Null
